# "Key Variables in Source Data" originally had an empty leading row (row 1)
# above its real header row (old row 2: Dataset | Variable | Type | Details).
# The edit removes that blank leading row, shifting every row below it up by
# one (new row 1 = old row 2, new row 35 = old row 36, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key Variables in Source Data")
$ws.Activate()

$ws.Rows.Item(1).Delete()

# Re-apply the existing sort on the "Medicare" block's Variable column so the
# sheet's <sortState>/<sortCondition> bookkeeping reflects the new row numbers
# (B16:B26 -> B15:B25) instead of staying stale.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B15:B25"))
$ws.Sort.SetRange($ws.Range("B15:B25"))
$ws.Sort.Apply()

# Restore frozen header row (was frozen through the old row 2, now through
# the new row 1) and the saved selection.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true

$ws.Range("A13").Select()
